# The "Analiza obiektowa - identyfikacja klas" bullet list (numId=4) drops the
# "UserSession" entry and the "Logger" entry; the remaining bullet items keep
# their original relative order (Transaction, CashDispenser, CashStorage,
# PINManager, BankDatabase, CardReader, Screen/Keypad/ReceiptPrinter).
#
# We find each target bullet paragraph by its distinctive leading text and
# delete the whole paragraph (including its paragraph mark) so that the
# following paragraphs shift up to take its place.

$d = $word.ActiveDocument

function Remove-ParagraphStartingWith([string]$prefix) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "$prefix*") {
            $p.Range.Delete()
            return $true
        }
    }
    return $false
}

Remove-ParagraphStartingWith("UserSession") | Out-Null
Remove-ParagraphStartingWith("Logger") | Out-Null
